{"js": "const replacements = [\n  [\"2025-09-12 Friday\", \"2025-09-13 Saturday\"],\n  [\"96\u00d766=6336\", \"37\u00d721=777\"],\n  [\"56\u00d738=2128\", \"41\u00d727=1107\"],\n  [\"44\u00d770=3080\", \"68\u00d793=6324\"],\n  [\"76\u00d742=3192\", \"29\u00d713=377\"],\n  [\"85\u00d772=6120\", \"62\u00d799=6138\"],\n  [\"37\u00d729=1073\", \"23\u00d751=1173\"],\n  [\"92\u00d714=1288\", \"66\u00d785=5610\"],\n  [\"69\u00d791=6279\", \"41\u00d729=1189\"],\n  [\"12\u00d766=792\", \"33\u00d765=2145\"],\n  [\"77\u00d724=1848\", \"52\u00d717=884\"],\n  [\"12\u00d761=732\", \"47\u00d736=1692\"],\n  [\"72\u00d726=1872\", \"68\u00d751=3468\"],\n  [\"60\u00d745=2700\", \"61\u00d733=2013\"],\n  [\"75\u00d732=2400\", \"24\u00d799=2376\"],\n  [\"73\u00d771=5183\", \"51\u00d778=3978\"],\n  [\"69\u00d753=3657\", \"40\u00d764=2560\"],\n  [\"77\u00d728=2156\", \"44\u00d762=2728\"],\n  [\"36\u00d774=2664\", \"17\u00d742=714\"],\n  [\"20\u00d738=760\", \"43\u00d799=4257\"],\n  [\"70\u00d752=3640\", \"64\u00d747=3008\"],\n  [\"20\u00d713=260\", \"18\u00d768=1224\"],\n  [\"96\u00d774=7104\", \"69\u00d732=2208\"],\n  [\"99\u00d724=2376\", \"41\u00d720=820\"],\n  [\"47\u00d719=893\", \"18\u00d735=630\"],\n  [\"70\u00d774=5180\", \"76\u00d735=2660\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-09-12 Friday\", \"2025-09-13 Saturday\"),\n  @(\"96\u00d766=6336\", \"37\u00d721=777\"),\n  @(\"56\u00d738=2128\", \"41\u00d727=1107\"),\n  @(\"44\u00d770=3080\", \"68\u00d793=6324\"),\n  @(\"76\u00d742=3192\", \"29\u00d713=377\"),\n  @(\"85\u00d772=6120\", \"62\u00d799=6138\"),\n  @(\"37\u00d729=1073\", \"23\u00d751=1173\"),\n  @(\"92\u00d714=1288\", \"66\u00d785=5610\"),\n  @(\"69\u00d791=6279\", \"41\u00d729=1189\"),\n  @(\"12\u00d766=792\", \"33\u00d765=2145\"),\n  @(\"77\u00d724=1848\", \"52\u00d717=884\"),\n  @(\"12\u00d761=732\", \"47\u00d736=1692\"),\n  @(\"72\u00d726=1872\", \"68\u00d751=3468\"),\n  @(\"60\u00d745=2700\", \"61\u00d733=2013\"),\n  @(\"75\u00d732=2400\", \"24\u00d799=2376\"),\n  @(\"73\u00d771=5183\", \"51\u00d778=3978\"),\n  @(\"69\u00d753=3657\", \"40\u00d764=2560\"),\n  @(\"77\u00d728=2156\", \"44\u00d762=2728\"),\n  @(\"36\u00d774=2664\", \"17\u00d742=714\"),\n  @(\"20\u00d738=760\", \"43\u00d799=4257\"),\n  @(\"70\u00d752=3640\", \"64\u00d747=3008\"),\n  @(\"20\u00d713=260\", \"18\u00d768=1224\"),\n  @(\"96\u00d774=7104\", \"69\u00d732=2208\"),\n  @(\"99\u00d724=2376\", \"41\u00d720=820\"),\n  @(\"47\u00d719=893\", \"18\u00d735=630\"),\n  @(\"70\u00d774=5180\", \"76\u00d735=2660\"),\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($pair[0], $false, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}"}
